$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G10").Value = "g"
$ws.Range("G11").Value = "h"
$ws.Range("G12").Value = "i"
$ws.Range("G13").Value = "j"
$ws.Range("G14").Value = "k"
$ws.Range("G15").Value = "l"
